$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 11.40722638115186

# Row 3
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 9.226618575922256
$ws.Range("D3").Value = 157.8057217802531
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 71685.65263718624

# Row 4
$ws.Range("B4").Value = 0.02258322285507441
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 9.924195562837738

# Row 5
$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 2938.103010863317
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 3188.247112057487

# Row 6
$ws.Range("B6").Value = 0.7287194209349384
$ws.Range("C6").Value = 9.226618575922256
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 246.9852506941017
$ws.Range("G6").Value = 260.0231881176624

# Row 7
$ws.Range("B7").Value = 0.3464964993005633
$ws.Range("C7").Value = 0.3375848360084654
$ws.Range("D7").Value = 3.082599426703578
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("G7").Value = 10.24810883928323

# Row 8
$ws.Range("B8").Value = 0.00006486019690155054
$ws.Range("C8").Value = 0.05231270169004087
$ws.Range("D8").Value = 116886.6739907443
$ws.Range("E8").Value = 71517.89157740913
$ws.Range("G8").Value = 188404.6179457154

# Row 9
$ws.Range("B9").Value = 0.000009318123435519965
$ws.Range("C9").Value = 0.3375848360084654
$ws.Range("D9").Value = 0.7127328510149897
$ws.Range("E9").Value = 5548678842208.939
$ws.Range("G9").Value = 5548678842209.99

# Row 10
$ws.Range("B10").Value = 0.02258322285507441
$ws.Range("C10").Value = 0.05231270169004087
$ws.Range("D10").Value = 157.8057217802531
$ws.Range("E10").Value = 71517.89157740913
$ws.Range("G10").Value = 71675.77219511392
